$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend the years table with columns M (2021) and N (2022) -----------
# Copy existing column L formatting into M and N so the new cells inherit
# the same cell styles (number formats, borders, fonts) as the rest of the
# table (rows 4-12).
$ws.Range("L4:L12").Copy()
$ws.Range("M4:N12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 4 - header years
$ws.Range("M4").Value = 2021
$ws.Range("N4").Value = 2022

# Row 5 - Kyrgyz Republic
$ws.Range("M5").Value = 5.6
$ws.Range("N5").Value = 6.3

# Row 6 - Batken oblast
$ws.Range("M6").Value = 0.8
$ws.Range("N6").Value = 0.8

# Row 7 - Djalal-Abad oblast
$ws.Range("M7").Value = 1.9
$ws.Range("N7").Value = 2.4

# Row 8 - Ysyk-Kul oblast
$ws.Range("M8").Value = 0.7
$ws.Range("N8").Value = 0.7

# Row 9 - Naryn oblast
$ws.Range("M9").Value = 0.7
$ws.Range("N9").Value = 0.8

# Row 10 - Osh oblast
$ws.Range("M10").Value = 0.9
$ws.Range("N10").Value = 1

# Row 11 - Talas oblast
$ws.Range("M11").Value = 0.3
$ws.Range("N11").Value = 0.2

# Row 12 - Chui oblast
$ws.Range("M12").Value = 0.2
$ws.Range("N12").Value = 0.4

# --- New footnote row 14 ---------------------------------------------------
$ws.Range("B13").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B14").Value = "По данным лесоустройства 2022 года Лесной службы при Министерстве чрезвычайных ситуаций КР"
$ws.Rows.Item(14).RowHeight = 34.5

# --- Tidy up the active selection -----------------------------------------
$ws.Range("A1").Select()
